$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = -20.588
$ws.Range("A7").Value = -21.286
$ws.Range("B7").Value = 6.312
$ws.Range("B15").Value = 5.042000000000001
$ws.Range("A16").Value = -22.043
$ws.Range("C16").Value = -13.27
$ws.Range("C19").Value = -12.192
$ws.Range("B21").Value = 8.644000000000002
$ws.Range("B22").Value = 6.709999999999999
$ws.Range("B23").Value = 7.359999999999999
$ws.Range("A28").Value = -21.767
$ws.Range("A29").Value = -21.586
$ws.Range("A32").Value = -21.724
$ws.Range("B34").Value = 7.492999999999999
$ws.Range("E34").Value = 16.751
$ws.Range("C36").Value = -12.669
$ws.Range("A40").Value = -20.646
$ws.Range("B43").Value = 6.049000000000001
$ws.Range("E43").Value = 17.283
$ws.Range("B45").Value = 5.672000000000001
$ws.Range("C46").Value = -14.113
$ws.Range("E48").Value = 17.205
$ws.Range("B50").Value = 5.715000000000001
$ws.Range("C50").Value = -13.493
$ws.Range("B51").Value = 6.284000000000001
$ws.Range("A52").Value = -21.63
$ws.Range("A57").Value = -22.32
$ws.Range("A66").Value = -21.322
$ws.Range("B66").Value = 5.83
$ws.Range("B67").Value = 5.571
$ws.Range("E70").Value = 17.387
$ws.Range("E73").Value = 16.485
$ws.Range("B79").Value = 5.502
$ws.Range("B84").Value = 5.453000000000001
$ws.Range("E87").Value = 16.465
$ws.Range("B92").Value = 5.247000000000001
$ws.Range("E92").Value = 17.916
$ws.Range("C95").Value = -12.18
$ws.Range("B97").Value = 6.811
$ws.Range("C97").Value = -13.606
$ws.Range("A100").Value = -21.907
$ws.Range("E101").Value = 16.569

$wb.Save()
